$wb = $excel.ActiveWorkbook

# Sheet "展览" — update "想去人数" (F column) counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value  = 2686
$ws1.Range("F10").Value = 122
$ws1.Range("F11").Value = 10129
$ws1.Range("F15").Value = 632
$ws1.Range("F16").Value = 11783
$ws1.Range("F17").Value = 12161
$ws1.Range("F19").Value = 96

# Sheet "全部类型" — same events, different row offsets
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value  = 2686
$ws4.Range("F11").Value = 122
$ws4.Range("F12").Value = 10129
$ws4.Range("F16").Value = 632
$ws4.Range("F17").Value = 11783
$ws4.Range("F18").Value = 12161
$ws4.Range("F20").Value = 96
